# "upgrade lvl calc": extend the per-day/per-night rate calc in Лист1 with
# two more derived rows (12-13), mirroring the existing pattern in rows 6-7
# but driven off hours (F) instead of days (E).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("I12").Formula = "=ROUND(H6/F6, 2)"
$ws.Range("J12").Formula = "=I12*F6"

$ws.Range("I13").Formula = "=ROUND(H7/F7, 2)"
$ws.Range("J13").Formula = "=I13*F7"

# Reflect the refreshed view state left by the edit (rezoom + reselect).
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$ws.Range("J6").Select()
